$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overlay")

# ---------------------------------------------------------------------------
# Add a new "Case Test" block of rows (18-23 in the '#' column, sheet rows
# 19-24) that demonstrate the new changeTextCase() function (lower / UPPER /
# Title). Row 19's formatting (A=number style, B:E=text style) is identical
# to the style already used by row 5, so clone it with Copy before writing
# the real values - this keeps the same shared cellXfs Excel already has
# instead of minting new ones.
# ---------------------------------------------------------------------------

$srcRow = $ws.Range("A5:E5")
for ($r = 19; $r -le 24; $r++) {
    $destRow = $ws.Range("A$r`:E$r")
    $srcRow.Copy($destRow)
    $ws.Rows.Item($r).RowHeight = 15.4
}

# Row 19 - header/config row for the new "Case Test" entry
$ws.Cells.Item(19, 1).Value = 18
$ws.Cells.Item(19, 2).Value = "Case Test"
$ws.Cells.Item(19, 3).Value = "<Type=Text><Text=Case Change Text Lower = >"
$ws.Cells.Item(19, 4).Value = "<X=25><Y=390>"
$ws.Cells.Item(19, 5).Value = "<Function=AddSpace(text,1)>"

# Row 20 - lower-case conversion
$ws.Cells.Item(20, 1).Value = 19
$ws.Cells.Item(20, 2).Value = "!<CONCAT><Case Test>"
$ws.Cells.Item(20, 3).Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=C>"
$ws.Cells.Item(20, 4).Value = ""
$ws.Cells.Item(20, 5).Value = "<Function=changeTextCase(text,lower)>"

# Row 21 - ", Upper = " label
$ws.Cells.Item(21, 1).Value = 20
$ws.Cells.Item(21, 2).Value = "!<CONCAT><Case Test>"
$ws.Cells.Item(21, 3).Value = "<Type=Text><Text=, Upper = >"
$ws.Cells.Item(21, 4).Value = ""
$ws.Cells.Item(21, 5).Value = "<Function=AddSpace(text,1)>"

# Row 22 - upper-case conversion
$ws.Cells.Item(22, 1).Value = 21
$ws.Cells.Item(22, 2).Value = "!<CONCAT><Case Test>"
$ws.Cells.Item(22, 3).Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=C>"
$ws.Cells.Item(22, 4).Value = ""
$ws.Cells.Item(22, 5).Value = "<Function=changeTextCase(text,UPPER)>"

# Row 23 - ", Title = " label
$ws.Cells.Item(23, 1).Value = 22
$ws.Cells.Item(23, 2).Value = "!<CONCAT><Case Test>"
$ws.Cells.Item(23, 3).Value = "<Type=Text><Text=, Title = >"
$ws.Cells.Item(23, 4).Value = ""
$ws.Cells.Item(23, 5).Value = "<Function=AddSpace(text,1)>"

# Row 24 - title-case conversion
$ws.Cells.Item(24, 1).Value = 23
$ws.Cells.Item(24, 2).Value = "!<CONCAT><Case Test>"
$ws.Cells.Item(24, 3).Value = "<Type=File><File=PAY01.xlsx><Sheet=SALERY DATA><PrimeryKey=A><Value=C>"
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = "<Function=changeTextCase(text,Title)>"
